$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 2728.2173
$ws.Range("J17").Value = 2728.2173
$ws.Range("L17").Value = 8184.651899999999
$ws.Range("N17").Value = -8520.651899999999

# Row 92
$ws.Range("H92").Value = 47619588
$ws.Range("J92").Value = 159
$ws.Range("L92").Value = 159
$ws.Range("N92").Value = -2655

# Row 113
$ws.Range("H113").Value = 3668.8096
$ws.Range("I113").Value = 2881.125
$ws.Range("J113").Value = 4153.5386
$ws.Range("K113").Value = 2881.125
$ws.Range("L113").Value = 4153.5386
$ws.Range("M113").Value = 372.875
$ws.Range("N113").Value = -10661.5386

# Row 132
$ws.Range("H132").Value = 12435.746
$ws.Range("I132").Value = 4761.59
$ws.Range("K132").Value = 14284.77
$ws.Range("M132").Value = -11754.77

# Row 135
$ws.Range("H135").Value = 3274.1155
$ws.Range("I135").Value = 1194.9375
$ws.Range("J135").Value = 6600.8
$ws.Range("K135").Value = 10754.4375
$ws.Range("L135").Value = 59407.2
$ws.Range("M135").Value = -8219.4375
$ws.Range("N135").Value = -64477.2

# Row 137
$ws.Range("H137").Value = 16354653
$ws.Range("I137").Value = 2002179.6
$ws.Range("K137").Value = 6006538.800000001
$ws.Range("M137").Value = -6003988.800000001

# Row 138
$ws.Range("H138").Value = 2015.9756
$ws.Range("I138").Value = 1078.8334
$ws.Range("J138").Value = 2556.6345
$ws.Range("K138").Value = 3236.5002
$ws.Range("L138").Value = 7669.9035
$ws.Range("M138").Value = 1903.4998
$ws.Range("N138").Value = -17949.9035

# Row 139
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# Row 140
$ws.Range("H140").Value = 65993.22
$ws.Range("J140").Value = 65403.75
$ws.Range("L140").Value = 65403.75
$ws.Range("N140").Value = -75763.75

# Row 141
$ws.Range("H141").Value = 3876.524
$ws.Range("I141").Value = 3876.524
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 11629.572
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -6449.572
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 12456.52
$ws.Range("I32").Value = 12663.295
$ws.Range("J32").Value = 11668.1875
$ws.Range("K32").Value = 12663.295
$ws.Range("L32").Value = 11668.1875
$ws.Range("M32").Value = -12376.295
$ws.Range("N32").Value = -12242.1875

# Row 45
$ws.Range("H45").Value = 5696.077
$ws.Range("I45").Value = 5548.8335
$ws.Range("K45").Value = 5548.8335
$ws.Range("M45").Value = -5171.8335

# Row 110
$ws.Range("H110").Value = 1460184.4
$ws.Range("I110").Value = 1701131.9
$ws.Range("K110").Value = 1701131.9
$ws.Range("M110").Value = -1699086.9

# Row 132
$ws.Range("H132").Value = 3698.0938
$ws.Range("I132").Value = 1630.775
$ws.Range("K132").Value = 4892.325000000001
$ws.Range("M132").Value = -2362.325000000001

$ws = $wb.Worksheets.Item("BSM")
# Row 105
$ws.Range("H105").Value = 2700.2942
$ws.Range("I105").Value = 2544.5
$ws.Range("K105").Value = 2544.5
$ws.Range("M105").Value = -797.5

# Row 116
$ws.Range("H116").Value = 79000
$ws.Range("J116").Value = 79000
$ws.Range("L116").Value = 79000
$ws.Range("N116").Value = -88178

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 6340.2373
$ws.Range("I31").Value = 2128.2222
$ws.Range("J31").Value = 7098.4
$ws.Range("K31").Value = 2128.2222
$ws.Range("L31").Value = 7098.4
$ws.Range("M31").Value = -1833.2222
$ws.Range("N31").Value = -7688.4

# Row 34
$ws.Range("H34").Value = 6340.2373
$ws.Range("I34").Value = 2128.2222
$ws.Range("J34").Value = 7098.4
$ws.Range("K34").Value = 2128.2222
$ws.Range("L34").Value = 7098.4
$ws.Range("M34").Value = -1926.2222
$ws.Range("N34").Value = -7502.4

# Row 58
$ws.Range("H58").Value = 5571.3213
$ws.Range("I58").Value = 2593.9092
$ws.Range("K58").Value = 2593.9092
$ws.Range("M58").Value = -2390.9092

# Row 105
$ws.Range("H105").Value = 1749469.4
$ws.Range("I105").Value = 3247602.2
$ws.Range("K105").Value = 3247602.2
$ws.Range("M105").Value = -3245855.2

# Row 132
$ws.Range("H132").Value = 8556288
$ws.Range("I132").Value = 9531251
$ws.Range("J132").Value = 25367.25
$ws.Range("K132").Value = 28593753
$ws.Range("L132").Value = 76101.75
$ws.Range("M132").Value = -28591223
$ws.Range("N132").Value = -81161.75

# Row 134
$ws.Range("H134").Value = 2215.3794
$ws.Range("I134").Value = 1769.96
$ws.Range("J134").Value = 4999.25
$ws.Range("K134").Value = 5309.88
$ws.Range("L134").Value = 14997.75
$ws.Range("M134").Value = -2774.88
$ws.Range("N134").Value = -20067.75

# Row 136
$ws.Range("H136").Value = 5571.3213
$ws.Range("I136").Value = 2593.9092
$ws.Range("K136").Value = 7781.7276
$ws.Range("M136").Value = -5231.7276

# Row 141
$ws.Range("H141").Value = 80720.85000000001
$ws.Range("J141").Value = 84620.83
$ws.Range("L141").Value = 84620.83
$ws.Range("N141").Value = -94980.83

$ws = $wb.Worksheets.Item("CUL")
# Row 48
$ws.Range("H48").Value = 5955.1113
$ws.Range("I48").Value = 900
$ws.Range("K48").Value = 2700
$ws.Range("M48").Value = -2450

# Row 113
$ws.Range("H113").Value = 768.125
$ws.Range("I113").Value = 704.8
$ws.Range("K113").Value = 2114.4
$ws.Range("M113").Value = 55.60000000000036

$ws = $wb.Worksheets.Item("GSM")
# Row 132
$ws.Range("H132").Value = 62687.44
$ws.Range("I132").Value = 76717.664
$ws.Range("J132").Value = 8570.857
$ws.Range("K132").Value = 230152.992
$ws.Range("L132").Value = 25712.571
$ws.Range("M132").Value = -227622.992
$ws.Range("N132").Value = -30772.571

$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 4060.7908
$ws.Range("I132").Value = 2846.9644
$ws.Range("J132").Value = 6326.6
$ws.Range("K132").Value = 8540.893199999999
$ws.Range("L132").Value = 18979.8
$ws.Range("M132").Value = -6010.893199999999
$ws.Range("N132").Value = -24039.8

# Row 135
$ws.Range("H135").Value = 39429
$ws.Range("J135").Value = 39429
$ws.Range("L135").Value = 39429
$ws.Range("N135").Value = -49569

# Row 136
$ws.Range("H136").Value = 3062.5
$ws.Range("I136").Value = 1534.091
$ws.Range("J136").Value = 8666.666999999999
$ws.Range("K136").Value = 4602.272999999999
$ws.Range("L136").Value = 26000.001
$ws.Range("M136").Value = -2052.272999999999
$ws.Range("N136").Value = -31100.001

# Row 137
$ws.Range("H137").Value = 22091.95
$ws.Range("J137").Value = 60429
$ws.Range("L137").Value = 60429
$ws.Range("N137").Value = -70629

$ws = $wb.Worksheets.Item("WVR")
# Row 107
$ws.Range("H107").Value = 3517.0356
$ws.Range("I107").Value = 4104.278
$ws.Range("K107").Value = 12312.834
$ws.Range("M107").Value = -10392.834

# Row 122
$ws.Range("H122").Value = 5166
$ws.Range("J122").Value = 6830.2
$ws.Range("L122").Value = 20490.6
$ws.Range("N122").Value = -25390.6

# Row 126
$ws.Range("H126").Value = 2735.0454
$ws.Range("I126").Value = 1401.0769
$ws.Range("K126").Value = 4203.2307
$ws.Range("M126").Value = -1733.2307

# Row 132
$ws.Range("H132").Value = 12824646
$ws.Range("I132").Value = 1409
$ws.Range("K132").Value = 4227
$ws.Range("M132").Value = -1697

# Row 136
$ws.Range("H136").Value = 8852.779
$ws.Range("I136").Value = 3017.2273
$ws.Range("J136").Value = 11187
$ws.Range("K136").Value = 9051.6819
$ws.Range("L136").Value = 33561
$ws.Range("M136").Value = -6501.6819
$ws.Range("N136").Value = -38661
